# Auto-generated edit script: updates cryptos list values to match the
# Fri Jul  5 14:56:29 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force every assigned string to be stored as literal text
# (prevents Excel from auto-converting numeric-looking strings like
# "1.00" or "492.02" into numbers, matching the original inlineStr cells).
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

Set-TextValue $ws.Range("D2") "55.780.08"
Set-TextValue $ws.Range("E2") "  -2.96%  "
Set-TextValue $ws.Range("D3") "2.958.52"
Set-TextValue $ws.Range("E3") "  -5.34%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "492.02"
Set-TextValue $ws.Range("E5") "  -5.63%  "
Set-TextValue $ws.Range("D6") "133.05"
Set-TextValue $ws.Range("E6") "  -0.86%  "
Set-TextValue $ws.Range("E7") "  -0.15%  "
Set-TextValue $ws.Range("D8") "2.956.50"
Set-TextValue $ws.Range("E8") "  -5.37%  "
Set-TextValue $ws.Range("E9") "  -5.42%  "
Set-TextValue $ws.Range("E10") "  -0.51%  "
Set-TextValue $ws.Range("D11") "0.102"
Set-TextValue $ws.Range("E11") "  -6.52%  "
Set-TextValue $ws.Range("E12") "  -9.26%  "
Set-TextValue $ws.Range("E13") "  +0.31%  "
Set-TextValue $ws.Range("D14") "3.464.63"
Set-TextValue $ws.Range("E14") "  -5.16%  "
Set-TextValue $ws.Range("D15") "24.75"
Set-TextValue $ws.Range("E15") "  -2.38%  "
Set-TextValue $ws.Range("D16") "55.732.58"
Set-TextValue $ws.Range("E16") "  -2.95%  "
Set-TextValue $ws.Range("D17") "2.957.88"
Set-TextValue $ws.Range("E17") "  -5.05%  "
Set-TextValue $ws.Range("E18") "  -5.57%  "
Set-TextValue $ws.Range("D19") "5.70"
Set-TextValue $ws.Range("E19") "  -1.09%  "
Set-TextValue $ws.Range("D20") "12.17"
Set-TextValue $ws.Range("E20") "  -5.85%  "
Set-TextValue $ws.Range("D21") "7.59"
Set-TextValue $ws.Range("E21") "  -4.91%  "
Set-TextValue $ws.Range("D22") "317.90"
Set-TextValue $ws.Range("E22") "  -7.37%  "
Set-TextValue $ws.Range("E23") "  +0.03%  "
Set-TextValue $ws.Range("E24") "  -8.52%  "
Set-TextValue $ws.Range("D25") "60.32"
Set-TextValue $ws.Range("E25") "  -11.66%  "
Set-TextValue $ws.Range("D26") "1.01"
Set-TextValue $ws.Range("E26") "  +0.58%  "
Set-TextValue $ws.Range("D27") "0.161"
Set-TextValue $ws.Range("E27") "  -2.96%  "
Set-TextValue $ws.Range("E28") "  -0.01%  "
Set-TextValue $ws.Range("D29") "0.0₃0847"
Set-TextValue $ws.Range("E29") "  -9.28%  "
Set-TextValue $ws.Range("D30") "6.48"
Set-TextValue $ws.Range("E30") "  -3.87%  "
Set-TextValue $ws.Range("D31") "6.61"
Set-TextValue $ws.Range("E31") "  -5.07%  "
Set-TextValue $ws.Range("D32") "1.15"
Set-TextValue $ws.Range("E32") "  -6.54%  "
Set-TextValue $ws.Range("E33") "  -8.41%  "
Set-TextValue $ws.Range("D34") "19.52"
Set-TextValue $ws.Range("E34") "  -8.98%  "
Set-TextValue $ws.Range("D35") "149.72"
Set-TextValue $ws.Range("E35") "  -5.66%  "
Set-TextValue $ws.Range("D36") "4.39"
Set-TextValue $ws.Range("E36") "  -8.11%  "
Set-TextValue $ws.Range("E37") "  -5.19%  "
Set-TextValue $ws.Range("E38") "  -7.82%  "
Set-TextValue $ws.Range("B39") "EnergySwap"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D39") "23.46"
Set-TextValue $ws.Range("E39") "  -6.74%  "
Set-TextValue $ws.Range("B40") "Hedera"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.0656"
Set-TextValue $ws.Range("E40") "  -4.46%  "
Set-TextValue $ws.Range("D41") "2.990.33"
Set-TextValue $ws.Range("E41") "  -5.02%  "
Set-TextValue $ws.Range("D42") "36.44"
Set-TextValue $ws.Range("E42") "  -9.55%  "
Set-TextValue $ws.Range("E43") "  -0.08%  "
Set-TextValue $ws.Range("B44") "Mantle"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D44") "0.631"
Set-TextValue $ws.Range("E44") "  -7.27%  "
Set-TextValue $ws.Range("B45") "ONDO"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D45") "0.989"
Set-TextValue $ws.Range("E45") "  -7.09%  "
Set-TextValue $ws.Range("D46") "1.39"
Set-TextValue $ws.Range("E46") "  -3.83%  "
Set-TextValue $ws.Range("E47") "  -9.04%  "
Set-TextValue $ws.Range("D48") "2.116.65"
Set-TextValue $ws.Range("E49") "  +0.70%  "
Set-TextValue $ws.Range("E50") "  -2.69%  "
Set-TextValue $ws.Range("B51") "dogwifhat"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D51") "1.83"
Set-TextValue $ws.Range("E51") "  +3.11%  "
